$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.6113651253405055
$ws.Range("J2").Value = 0.6113651253405055
$ws.Range("M2").Value = 1.620350333333333
$ws.Range("N2").Value = 4.861051
$ws.Range("O2").Value = 0.0725197794467048
$ws.Range("P2").Value = 0.07251977944670479
$ws.Range("Q2").Value = 0.04788189246677777
$ws.Range("R2").Value = 0.4309370322009999
$ws.Range("S2").Value = 0.04433606405110049
$ws.Range("T2").Value = 0.04433606405110049
$ws.Range("I3").Value = 0.6113651253405055
$ws.Range("J3").Value = 0.6113651253405055
$ws.Range("O3").Value = 0.09715752300453648
$ws.Range("P3").Value = 0.09715752300453648
$ws.Range("S3").Value = 0.05939872122944149
$ws.Range("T3").Value = 0.05939872122944149
$ws.Range("I4").Value = 0.6113651253405055
$ws.Range("J4").Value = 0.6113651253405055
$ws.Range("M4").Value = 13.65106133333333
$ws.Range("N4").Value = 40.953184
$ws.Range("O4").Value = 0.6109616770777183
$ws.Range("P4").Value = 0.6109616770777183
$ws.Range("Q4").Value = 0.4033934127537777
$ws.Range("R4").Value = 3.630540714784
$ws.Range("S4").Value = 0.3735206622848647
$ws.Range("T4").Value = 0.3735206622848647
$ws.Range("I5").Value = 0.6113651253405055
$ws.Range("J5").Value = 0.6113651253405055
$ws.Range("M5").Value = 0.6533493333333333
$ws.Range("N5").Value = 1.960048
$ws.Range("O5").Value = 0.02924105274043717
$ws.Range("P5").Value = 0.02924105274043717
$ws.Range("Q5").Value = 0.01930669058311111
$ws.Range("R5").Value = 0.173760215248
$ws.Range("S5").Value = 0.0178769598737457
$ws.Range("T5").Value = 0.0178769598737457
$ws.Range("I6").Value = 0.6113651253405055
$ws.Range("J6").Value = 0.6113651253405055
$ws.Range("M6").Value = 4.247957666666667
$ws.Range("N6").Value = 12.743873
$ws.Range("O6").Value = 0.1901199677306032
$ws.Range("P6").Value = 0.1901199677306032
$ws.Range("Q6").Value = 0.1255285650358889
$ws.Range("R6").Value = 1.129757085323
$ws.Range("S6").Value = 0.1162327179013531
$ws.Range("T6").Value = 0.1162327179013531
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.01878466666666667
$ws.Range("H7").Value = 0.056354
$ws.Range("I7").Value = 0.3886348746594945
$ws.Range("J7").Value = 0.3886348746594945
$ws.Range("M7").Value = 1.620350333333333
$ws.Range("N7").Value = 4.861051
$ws.Range("O7").Value = 0.0725197794467048
$ws.Range("P7").Value = 0.07251977944670479
$ws.Range("Q7").Value = 0.03043774089488889
$ws.Range("R7").Value = 0.273939668054
$ws.Range("S7").Value = 0.02818371539560431
$ws.Range("T7").Value = 0.0281837153956043
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.01878466666666667
$ws.Range("H8").Value = 0.056354
$ws.Range("I8").Value = 0.3886348746594945
$ws.Range("J8").Value = 0.3886348746594945
$ws.Range("O8").Value = 0.09715752300453648
$ws.Range("P8").Value = 0.09715752300453648
$ws.Range("Q8").Value = 0.04077860597155556
$ws.Range("R8").Value = 0.367007453744
$ws.Range("S8").Value = 0.03775880177509499
$ws.Range("T8").Value = 0.03775880177509499
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.01878466666666667
$ws.Range("H9").Value = 0.056354
$ws.Range("I9").Value = 0.3886348746594945
$ws.Range("J9").Value = 0.3886348746594945
$ws.Range("M9").Value = 13.65106133333333
$ws.Range("N9").Value = 40.953184
$ws.Range("O9").Value = 0.6109616770777183
$ws.Range("P9").Value = 0.6109616770777183
$ws.Range("Q9").Value = 0.2564306367928889
$ws.Range("R9").Value = 2.307875731136
$ws.Range("S9").Value = 0.2374410147928536
$ws.Range("T9").Value = 0.2374410147928536
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.01878466666666667
$ws.Range("H10").Value = 0.056354
$ws.Range("I10").Value = 0.3886348746594945
$ws.Range("J10").Value = 0.3886348746594945
$ws.Range("M10").Value = 0.6533493333333333
$ws.Range("N10").Value = 1.960048
$ws.Range("O10").Value = 0.02924105274043717
$ws.Range("P10").Value = 0.02924105274043717
$ws.Range("Q10").Value = 0.01227294944355556
$ws.Range("R10").Value = 0.110456544992
$ws.Range("S10").Value = 0.01136409286669147
$ws.Range("T10").Value = 0.01136409286669147
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.01878466666666667
$ws.Range("H11").Value = 0.056354
$ws.Range("I11").Value = 0.3886348746594945
$ws.Range("J11").Value = 0.3886348746594945
$ws.Range("M11").Value = 4.247957666666667
$ws.Range("N11").Value = 12.743873
$ws.Range("O11").Value = 0.1901199677306032
$ws.Range("P11").Value = 0.1901199677306032
$ws.Range("Q11").Value = 0.07979646878244447
$ws.Range("R11").Value = 0.718168219042
$ws.Range("S11").Value = 0.07388724982925012
$ws.Range("T11").Value = 0.0738872498292501
